$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Bad Drivers" table values per weekly driver report refresh
$ws.Range("C3").Value = 2537
$ws.Range("D3").Value = 94.90000000000001
$ws.Range("C4").Value = 501
$ws.Range("C5").Value = 3038
